# Apply narrative cleanup to the FM-multiBirth StructureDefinition workbook.
# Renames "multiBirth" -> "MultiBirth" in the URL/Name/Title/Description text,
# updates the "multi-birth" casing in the title/description, and bumps the
# publication Date to the new commit timestamp.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet (Property / Value pairs) ---

# URL
$metadata.Range("B2").Value = "https://johnmoehrke.github.io/MHV-PGHD/StructureDefinition/FM-MultiBirth"

# Name
$metadata.Range("B4").Value = "MultiBirth"

# Title
$metadata.Range("B5").Value = "MultiBirth indication"

# Date
$metadata.Range("B8").Value = "2022-04-11T07:37:02-05:00"

# Description
$metadata.Range("B12").Value = "When this family member is known to be part of a MultiBirth, indicate how many siblings."

# --- Elements sheet (same shared strings appear in the Short/Definition/Fixed Value columns) ---

# Short (row 2, column K) mirrors the Title text
$elements.Range("K2").Value = "MultiBirth indication"

# Definition (row 2, column L) mirrors the Description text
$elements.Range("L2").Value = "When this family member is known to be part of a MultiBirth, indicate how many siblings."

# Fixed Value (row 5, column Q) mirrors the URL text
$elements.Range("Q5").Value = "https://johnmoehrke.github.io/MHV-PGHD/StructureDefinition/FM-MultiBirth"
